$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1697792869269949
$ws.Range("C2").Value = 0.599320882852292
$ws.Range("J2").Value = 0.01188455008488964
$ws.Range("P2").Value = 0.1290322580645161
$ws.Range("S2").Value = 0.0899830220713073
$ws.Range("B3").Value = 0.01063829787234043
$ws.Range("C3").Value = 0.05053191489361702
$ws.Range("J3").Value = 0.02127659574468085
$ws.Range("P3").Value = 0.6941489361702128
$ws.Range("S3").Value = 0.2234042553191489
$ws.Range("J4").Value = 0.06930693069306931
$ws.Range("P4").Value = 0.7029702970297029
$ws.Range("S4").Value = 0.2277227722772277
$ws.Range("P5").Value = 0.25
$ws.Range("S5").Value = 0.75
$ws.Range("B6").Value = 0.064453125
$ws.Range("D6").Value = 0.01953125
$ws.Range("E6").Value = 0.00390625
$ws.Range("F6").Value = 0.08203125
$ws.Range("J6").Value = 0.24609375
$ws.Range("O6").Value = 0.017578125
$ws.Range("Q6").Value = 0.15625
$ws.Range("R6").Value = 0.064453125
$ws.Range("S6").Value = 0.345703125
$ws.Range("B7").Value = 0.1132478632478632
$ws.Range("D7").Value = 0.01495726495726496
$ws.Range("F7").Value = 0.04700854700854701
$ws.Range("J7").Value = 0.1474358974358974
$ws.Range("O7").Value = 0.01282051282051282
$ws.Range("Q7").Value = 0.1645299145299145
$ws.Range("R7").Value = 0.07692307692307693
$ws.Range("S7").Value = 0.4230769230769231
$ws.Range("B8").Value = 0.08
$ws.Range("D8").Value = 0.01333333333333333
$ws.Range("E8").Value = 0.0009523809523809524
$ws.Range("F8").Value = 0.05047619047619047
$ws.Range("J8").Value = 0.1123809523809524
$ws.Range("O8").Value = 0.02380952380952381
$ws.Range("Q8").Value = 0.1828571428571429
$ws.Range("R8").Value = 0.08761904761904762
$ws.Range("S8").Value = 0.4485714285714286
$ws.Range("B9").Value = 0.07658643326039387
$ws.Range("D9").Value = 0.01531728665207877
$ws.Range("F9").Value = 0.05908096280087528
$ws.Range("J9").Value = 0.137855579868709
$ws.Range("O9").Value = 0.02407002188183808
$ws.Range("Q9").Value = 0.2013129102844639
$ws.Range("R9").Value = 0.09409190371991247
$ws.Range("S9").Value = 0.3916849015317287
$ws.Range("B10").Value = 0.0986013986013986
$ws.Range("D10").Value = 0.02377622377622378
$ws.Range("E10").Value = 0.001048951048951049
$ws.Range("F10").Value = 0.07727272727272727
$ws.Range("J10").Value = 0.1251748251748252
$ws.Range("O10").Value = 0.01608391608391608
$ws.Range("Q10").Value = 0.2066433566433566
$ws.Range("R10").Value = 0.07902097902097902
$ws.Range("S10").Value = 0.3723776223776224
$ws.Range("G11").Value = 0.1518105849582173
$ws.Range("J11").Value = 0.08356545961002786
$ws.Range("K11").Value = 0.2047353760445682
$ws.Range("L11").Value = 0.5431754874651811
$ws.Range("S11").Value = 0.01671309192200557
$ws.Range("G12").Value = 0.7518796992481203
$ws.Range("J12").Value = 0.1854636591478697
$ws.Range("K12").Value = 0.007518796992481203
$ws.Range("L12").Value = 0.01503759398496241
$ws.Range("S12").Value = 0.04010025062656641
$ws.Range("G13").Value = 0.6576576576576577
$ws.Range("J13").Value = 0.3153153153153153
$ws.Range("S13").Value = 0.02702702702702703
$ws.Range("G14").Value = 0.6
$ws.Range("J14").Value = 0.4
$ws.Range("F15").Value = 0.01232032854209446
$ws.Range("H15").Value = 0.1581108829568789
$ws.Range("I15").Value = 0.08213552361396304
$ws.Range("J15").Value = 0.3305954825462012
$ws.Range("K15").Value = 0.08418891170431211
$ws.Range("M15").Value = 0.002053388090349076
$ws.Range("O15").Value = 0.06160164271047228
$ws.Range("S15").Value = 0.2689938398357289
$ws.Range("F16").Value = 0.02010050251256281
$ws.Range("H16").Value = 0.1733668341708543
$ws.Range("I16").Value = 0.09296482412060302
$ws.Range("J16").Value = 0.3944723618090452
$ws.Range("K16").Value = 0.1080402010050251
$ws.Range("M16").Value = 0.02010050251256281
$ws.Range("O16").Value = 0.05025125628140704
$ws.Range("S16").Value = 0.1407035175879397
$ws.Range("F17").Value = 0.01958863858961802
$ws.Range("H17").Value = 0.1890303623898139
$ws.Range("I17").Value = 0.09696376101860921
$ws.Range("J17").Value = 0.4035259549461312
$ws.Range("K17").Value = 0.09990205680705191
$ws.Range("M17").Value = 0.01469147894221352
$ws.Range("N17").Value = 0.0009794319294809011
$ws.Range("O17").Value = 0.05876591576885407
$ws.Range("S17").Value = 0.1165523996082272
$ws.Range("F18").Value = 0.02790697674418605
$ws.Range("H18").Value = 0.1906976744186047
$ws.Range("I18").Value = 0.09302325581395349
$ws.Range("J18").Value = 0.3953488372093023
$ws.Range("K18").Value = 0.1093023255813954
$ws.Range("M18").Value = 0.02325581395348837
$ws.Range("N18").Value = 0.004651162790697674
$ws.Range("O18").Value = 0.05116279069767442
$ws.Range("S18").Value = 0.1046511627906977
$ws.Range("F19").Value = 0.0161512027491409
$ws.Range("H19").Value = 0.2195876288659794
$ws.Range("I19").Value = 0.0831615120274914
$ws.Range("J19").Value = 0.3680412371134021
$ws.Range("K19").Value = 0.1130584192439863
$ws.Range("M19").Value = 0.02577319587628866
$ws.Range("N19").Value = 0.001030927835051546
$ws.Range("O19").Value = 0.06769759450171821
$ws.Range("S19").Value = 0.1054982817869416
